$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.060.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.298.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.656.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.345.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.944.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.20%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0696"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.008.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.51%  "
